$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 onto the new header cells I1:J1 so they match
# the existing bold/centered/bordered header formatting.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-57
$data = @(
    @{ Row = 2; I = 8; J = 8 },
    @{ Row = 3; I = 5; J = 6 },
    @{ Row = 4; I = 6; J = 6 },
    @{ Row = 5; I = 7; J = 7 },
    @{ Row = 6; I = 10; J = 10 },
    @{ Row = 7; I = 8; J = 8 },
    @{ Row = 8; I = 7; J = 7 },
    @{ Row = 9; I = 6; J = 7 },
    @{ Row = 10; I = 6; J = 6 },
    @{ Row = 11; I = 9; J = 9 },
    @{ Row = 12; I = 4; J = 5 },
    @{ Row = 13; I = 7; J = 7 },
    @{ Row = 14; I = 8; J = 8 },
    @{ Row = 15; I = 7; J = 7 },
    @{ Row = 16; I = 6; J = 6 },
    @{ Row = 17; I = 5; J = 6 },
    @{ Row = 18; I = 7; J = 7 },
    @{ Row = 19; I = 5; J = 5 },
    @{ Row = 20; I = 8; J = 8 },
    @{ Row = 21; I = 7; J = 7 },
    @{ Row = 22; I = 6; J = 7 },
    @{ Row = 23; I = 6; J = 7 },
    @{ Row = 24; I = 8; J = 8 },
    @{ Row = 25; I = 4; J = 5 },
    @{ Row = 26; I = 7; J = 8 },
    @{ Row = 27; I = 6; J = 7 },
    @{ Row = 28; I = 3; J = 4 },
    @{ Row = 29; I = 8; J = 8 },
    @{ Row = 30; I = 5; J = 6 },
    @{ Row = 31; I = 8; J = 8 },
    @{ Row = 32; I = 6; J = 7 },
    @{ Row = 33; I = 5; J = 6 },
    @{ Row = 34; I = 9; J = 9 },
    @{ Row = 35; I = 6; J = 7 },
    @{ Row = 36; I = 7; J = 8 },
    @{ Row = 37; I = 7; J = 8 },
    @{ Row = 38; I = 8; J = 8 },
    @{ Row = 39; I = 7; J = 7 },
    @{ Row = 40; I = 6; J = 6 },
    @{ Row = 41; I = 6; J = 7 },
    @{ Row = 42; I = 6; J = 7 },
    @{ Row = 43; I = 5; J = 7 },
    @{ Row = 44; I = 9; J = 9 },
    @{ Row = 45; I = 10; J = 10 },
    @{ Row = 46; I = 8; J = 8 },
    @{ Row = 47; I = 5; J = 6 },
    @{ Row = 48; I = 5; J = 6 },
    @{ Row = 49; I = 5; J = 7 },
    @{ Row = 50; I = 8; J = 9 },
    @{ Row = 51; I = 10; J = 10 },
    @{ Row = 52; I = 5; J = 5 },
    @{ Row = 53; I = 2; J = 3 },
    @{ Row = 54; I = 1; J = 5 },
    @{ Row = 55; I = 1; J = 4 },
    @{ Row = 56; I = 1; J = 3 },
    @{ Row = 57; I = 1; J = 2 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}
